# Revert "adding term 2.0 now utf-8"
# This reverts the workbook back to the 1.1.0 state:
#  - removes the 5 worksheets that were added for term 2.0
#    ("Include from FSIII 6" .. "Include from FSIII 10")
#  - restores the "Value" cell on the remaining "Include from FSIII*"
#    sheets from the FSIII UUIDs back to the plain J1..J5 codes
#  - restores the Metadata sheet's Version / Date / Contact values

$wb = $excel.ActiveWorkbook

# --- remove the sheets added by the 2.0 update ---------------------------
$sheetsToRemove = @(
    "Include from FSIII 6",
    "Include from FSIII 7",
    "Include from FSIII 8",
    "Include from FSIII 9",
    "Include from FSIII 10"
)
foreach ($name in $sheetsToRemove) {
    $wb.Worksheets.Item($name).Delete() | Out-Null
}

# --- restore the "Value" column on the remaining include sheets ----------
$valueBySheet = @{
    "Include from FSIII"   = "J1"
    "Include from FSIII 2" = "J2"
    "Include from FSIII 3" = "J3"
    "Include from FSIII 4" = "J4"
    "Include from FSIII 5" = "J5"
}
foreach ($name in $valueBySheet.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C2").Value = $valueBySheet[$name]
}

# --- restore the Metadata values ------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# keep the original active sheet/tab (Metadata, first tab) selected
$meta.Activate()
